$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 401757.2
$ws.Range("J17").Value = 401757.2
$ws.Range("L17").Value = 1205271.6
$ws.Range("N17").Value = -1205607.6

# Row 96
$ws.Range("H96").Value = 719.4
$ws.Range("I96").Value = 699.4286
$ws.Range("J96").Value = 999
$ws.Range("K96").Value = 2098.2858
$ws.Range("L96").Value = 2997
$ws.Range("M96").Value = -725.2857999999997
$ws.Range("N96").Value = -5743

# Row 99
$ws.Range("H99").Value = 66667108
$ws.Range("I99").Value = 499.2
$ws.Range("J99").Value = 200000340
$ws.Range("K99").Value = 1497.6
$ws.Range("L99").Value = 600001020
$ws.Range("M99").Value = 0.4000000000000909
$ws.Range("N99").Value = -600004016

# Row 101
$ws.Range("H101").Value = 935.7222
$ws.Range("I101").Value = 940.1875
$ws.Range("J101").Value = 900
$ws.Range("K101").Value = 2820.5625
$ws.Range("L101").Value = 2700
$ws.Range("M101").Value = -1198.5625
$ws.Range("N101").Value = -5944

# Row 116
$ws.Range("H116").Value = 4318.0454
$ws.Range("I116").Value = 3866.4
$ws.Range("K116").Value = 3866.4
$ws.Range("M116").Value = -424.4000000000001

# Row 137
$ws.Range("H137").Value = 948.55554
$ws.Range("I137").Value = 867.75
$ws.Range("J137").Value = 1179.4286
$ws.Range("K137").Value = 2603.25
$ws.Range("L137").Value = 3538.2858
$ws.Range("M137").Value = -53.25
$ws.Range("N137").Value = -8638.2858

# Row 138
$ws.Range("H138").Value = 2254.5312
$ws.Range("J138").Value = 2554.5625
$ws.Range("L138").Value = 7663.6875
$ws.Range("N138").Value = -17943.6875

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 10734.934
$ws.Range("I45").Value = 19017.834
$ws.Range("K45").Value = 19017.834
$ws.Range("M45").Value = -18640.834

# Row 61
$ws.Range("H61").Value = 2055.2144
$ws.Range("I61").Value = 1991.5555
$ws.Range("J61").Value = 2169.8
$ws.Range("K61").Value = 1991.5555
$ws.Range("L61").Value = 2169.8
$ws.Range("M61").Value = -1779.5555
$ws.Range("N61").Value = -2593.8

# Row 102
$ws.Range("H102").Value = 1928.8
$ws.Range("I102").Value = 1918.4688
$ws.Range("K102").Value = 1918.4688
$ws.Range("M102").Value = -296.4688000000001

# Row 136
$ws.Range("H136").Value = 2055.2144
$ws.Range("I136").Value = 1991.5555
$ws.Range("J136").Value = 2169.8
$ws.Range("K136").Value = 5974.666499999999
$ws.Range("L136").Value = 6509.400000000001
$ws.Range("M136").Value = -3424.666499999999
$ws.Range("N136").Value = -11609.4

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 717.5
$ws.Range("I80").Value = 537.5714
$ws.Range("J80").Value = 832
$ws.Range("K80").Value = 537.5714
$ws.Range("L80").Value = 832
$ws.Range("M80").Value = 460.4286
$ws.Range("N80").Value = -2828

# Row 83
$ws.Range("H83").Value = 717.5
$ws.Range("I83").Value = 537.5714
$ws.Range("J83").Value = 832
$ws.Range("K83").Value = 2687.857
$ws.Range("L83").Value = 4160
$ws.Range("M83").Value = 2304.143
$ws.Range("N83").Value = -14144

# Row 94
$ws.Range("H94").Value = 1005.8
$ws.Range("I94").Value = 1084.3846
$ws.Range("K94").Value = 1084.3846
$ws.Range("M94").Value = -633.3846000000001

# Row 99
$ws.Range("H99").Value = 2902.375
$ws.Range("I99").Value = 2036.5
$ws.Range("J99").Value = 5500
$ws.Range("K99").Value = 2036.5
$ws.Range("L99").Value = 5500
$ws.Range("M99").Value = -538.5
$ws.Range("N99").Value = -8496

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 840
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -88
$ws.Range("N4").Value = -1224

# Row 16
$ws.Range("H16").Value = 2751.25
$ws.Range("I16").Value = 2868.5
$ws.Range("K16").Value = 2868.5
$ws.Range("M16").Value = -2581.5

# Row 31
$ws.Range("H31").Value = 9838.120000000001
$ws.Range("I31").Value = 2926.861
$ws.Range("K31").Value = 2926.861
$ws.Range("M31").Value = -2631.861

# Row 34
$ws.Range("H34").Value = 9838.120000000001
$ws.Range("I34").Value = 2926.861
$ws.Range("K34").Value = 2926.861
$ws.Range("M34").Value = -2724.861

# Row 62
$ws.Range("H62").Value = 3356.7144
$ws.Range("J62").Value = 3374.75
$ws.Range("L62").Value = 3374.75
$ws.Range("N62").Value = -4622.75

# Row 65
$ws.Range("H65").Value = 3356.7144
$ws.Range("J65").Value = 3374.75
$ws.Range("L65").Value = 16873.75
$ws.Range("N65").Value = -23113.75

# Row 99
$ws.Range("H99").Value = 40599.855
$ws.Range("I99").Value = 67449.75
$ws.Range("K99").Value = 67449.75
$ws.Range("M99").Value = -65951.75

# Row 105
$ws.Range("H105").Value = 1810.2142
$ws.Range("I105").Value = 2124.3
$ws.Range("J105").Value = 1025
$ws.Range("K105").Value = 2124.3
$ws.Range("L105").Value = 1025
$ws.Range("M105").Value = -377.3000000000002
$ws.Range("N105").Value = -4519

# Row 113
$ws.Range("H113").Value = 2751.25
$ws.Range("I113").Value = 2868.5
$ws.Range("K113").Value = 2868.5
$ws.Range("M113").Value = -698.5

# Row 126
$ws.Range("H126").Value = 40599.855
$ws.Range("I126").Value = 67449.75
$ws.Range("K126").Value = 202349.25
$ws.Range("M126").Value = -199879.25

# Row 132
$ws.Range("H132").Value = 5135.294
$ws.Range("I132").Value = 5175.857
$ws.Range("J132").Value = 4946
$ws.Range("K132").Value = 15527.571
$ws.Range("L132").Value = 14838
$ws.Range("M132").Value = -12997.571
$ws.Range("N132").Value = -19898

# Row 141
$ws.Range("H141").Value = 375356.44
$ws.Range("J141").Value = 375356.44
$ws.Range("L141").Value = 375356.44
$ws.Range("N141").Value = -385716.44

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 819.9286
$ws.Range("I34").Value = 179.71428
$ws.Range("J34").Value = 1460.1428
$ws.Range("K34").Value = 539.14284
$ws.Range("L34").Value = 4380.428400000001
$ws.Range("M34").Value = -455.14284
$ws.Range("N34").Value = -4548.428400000001

# Row 39
$ws.Range("H39").Value = 4264.75
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4264.75
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 12794.25
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -13382.25

# Row 55
$ws.Range("H55").Value = 8335841
$ws.Range("I55").Value = 1199.75
$ws.Range("J55").Value = 11366620
$ws.Range("K55").Value = 3599.25
$ws.Range("L55").Value = 34099860
$ws.Range("M55").Value = -3422.25
$ws.Range("N55").Value = -34100214

# Row 114
$ws.Range("H114").Value = 50000724
$ws.Range("I114").Value = 50000724
$ws.Range("K114").Value = 150002172
$ws.Range("M114").Value = -149998918

# Row 129
$ws.Range("H129").Value = 156639.23
$ws.Range("J129").Value = 4380.4287
$ws.Range("L129").Value = 13141.2861
$ws.Range("N129").Value = -23141.2861

$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 17334.166

# Row 24
$ws.Range("H24").Value = 21175.45
$ws.Range("J24").Value = 21538.857
$ws.Range("L24").Value = 21538.857
$ws.Range("N24").Value = -21884.857

# Row 34
$ws.Range("H34").Value = 25173
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 25173
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 25173
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -25709

# Row 76
$ws.Range("H76").Value = 25173
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 25173
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 25173
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -25803

# Row 79
$ws.Range("H79").Value = 25173
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 25173
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 25173
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -27357

# Row 123
$ws.Range("H123").Value = 58743.273
$ws.Range("J123").Value = 58743.273
$ws.Range("L123").Value = 58743.273
$ws.Range("N123").Value = -63643.273

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 84563.914
$ws.Range("I16").Value = 1307.4445
$ws.Range("K16").Value = 1307.4445
$ws.Range("M16").Value = -1137.4445

# Row 22
$ws.Range("H22").Value = 695
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 690
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 690
$ws.Range("M22").Value = -405
$ws.Range("N22").Value = -1280

# Row 27
$ws.Range("H27").Value = 695
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 690
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 690
$ws.Range("M27").Value = -593
$ws.Range("N27").Value = -904

# Row 43
$ws.Range("H43").Value = 11753905
$ws.Range("J43").Value = 11753905
$ws.Range("L43").Value = 11753905
$ws.Range("N43").Value = -11754291

# Row 45
$ws.Range("H45").Value = 11669.167
$ws.Range("I45").Value = 2520
$ws.Range("K45").Value = 2520
$ws.Range("M45").Value = -2113

# Row 93
$ws.Range("H93").Value = 10303.185
$ws.Range("I93").Value = 1602.7941
$ws.Range("K93").Value = 1602.7941
$ws.Range("M93").Value = -354.7941000000001

# Row 107
$ws.Range("H107").Value = 3163
$ws.Range("I107").Value = 3163
$ws.Range("K107").Value = 3163
$ws.Range("M107").Value = -1243

# Row 136
$ws.Range("H136").Value = 4557.3335
$ws.Range("I136").Value = 3122.75
$ws.Range("K136").Value = 9368.25
$ws.Range("M136").Value = -6818.25

$ws = $wb.Worksheets.Item("WVR")
# Row 55
$ws.Range("H55").Value = 3500
$ws.Range("I55").Value = 3500
$ws.Range("K55").Value = 3500
$ws.Range("M55").Value = -3223

# Row 100
$ws.Range("H100").Value = 3078.1304
$ws.Range("I100").Value = 2815.5625
$ws.Range("K100").Value = 5631.125
$ws.Range("M100").Value = -5090.125

# Row 123
$ws.Range("H123").Value = 30551.666
$ws.Range("J123").Value = 30551.666
$ws.Range("L123").Value = 30551.666
$ws.Range("N123").Value = -40351.666

# Row 132
$ws.Range("H132").Value = 1914.3334
$ws.Range("I132").Value = 1407.8125
$ws.Range("J132").Value = 2493.2144
$ws.Range("K132").Value = 4223.4375
$ws.Range("L132").Value = 7479.6432
$ws.Range("M132").Value = -1693.4375
$ws.Range("N132").Value = -12539.6432
